{"js": "// Apply the text replacements described by the diff using the Word\n// JavaScript API (Office.js). `context` (alias `ctx`) is the\n// Word.RequestContext passed in by the harness.\n\nconst replacements = [\n  {\n    find: \"Play Frozen Queen for Free: A Unique Fantasy Slot\",\n    replace: \"Play Frozen Queen Free | Review of the Engaging Slot Game\",\n  },\n  {\n    find: \"Bonus features and symbols\",\n    replace: \"Generous bonus features and symbols\",\n  },\n  {\n    find: \"Exemplary fantasy-themed setting\",\n    replace: \"Immersive fantasy-themed setting\",\n  },\n  {\n    find: \"Only two ice crystal symbols required for activation\",\n    replace: \"Limited number of reels\",\n  },\n  {\n    find: \"Limited 84 pay lines\",\n    replace: \"May not appeal to players looking for traditional slot games\",\n  },\n  {\n    find:\n      \"Try Frozen Queen for free and experience the engaging gameplay mechanics, stunning graphics, and winter wonderland setting of this exemplary fantasy-themed slot game.\",\n    replace:\n      \"Read our review of Frozen Queen, an engaging and unique slot game to play for free.\",\n  },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = context.document.body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the text replacements described by the diff using the Word COM\n# object model. $word / $d resolve to the running application / the\n# open ActiveDocument.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = \"Play Frozen Queen for Free: A Unique Fantasy Slot\"; Replace = \"Play Frozen Queen Free | Review of the Engaging Slot Game\" },\n    @{ Find = \"Bonus features and symbols\"; Replace = \"Generous bonus features and symbols\" },\n    @{ Find = \"Exemplary fantasy-themed setting\"; Replace = \"Immersive fantasy-themed setting\" },\n    @{ Find = \"Only two ice crystal symbols required for activation\"; Replace = \"Limited number of reels\" },\n    @{ Find = \"Limited 84 pay lines\"; Replace = \"May not appeal to players looking for traditional slot games\" },\n    @{ Find = \"Try Frozen Queen for free and experience the engaging gameplay mechanics, stunning graphics, and winter wonderland setting of this exemplary fantasy-themed slot game.\"; Replace = \"Read our review of Frozen Queen, an engaging and unique slot game to play for free.\" }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $null = $find.Execute($r.Find, $true, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2)\n}\n\nWrite-Output \"done\"\n"}
